$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text storage for price cells whose updated values would
# otherwise be auto-recognized as numbers by Excel (so trailing zeros,
# e.g. "1.00" or "69.70", survive exactly like the rest of column D).
$textCells = @("D4", "D5", "D6", "D8", "D9", "D12", "D13", "D17", "D19", "D20", "D21", "D25", "D27", "D29", "D32", "D34", "D35", "D37", "D38", "D42", "D45", "D46", "D47")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "74.526.05"
$ws.Range("E2").Value = "  +8.75%  "
$ws.Range("D3").Value = "2.589.10"
$ws.Range("E3").Value = "  +6.60%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "184.56"
$ws.Range("E5").Value = "  +14.95%  "
$ws.Range("D6").Value = "579.43"
$ws.Range("E6").Value = "  +3.73%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").Value = "  +4.74%  "
$ws.Range("D9").Value = "0.204"
$ws.Range("E9").Value = "  +24.97%  "
$ws.Range("D10").Value = "2.588.10"
$ws.Range("E10").Value = "  +6.61%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "0.359"
$ws.Range("E12").Value = "  +8.64%  "
$ws.Range("D13").Value = "4.76"
$ws.Range("E13").Value = "  +3.39%  "
$ws.Range("E14").Value = "  +9.82%  "
$ws.Range("D15").Value = "74.358.42"
$ws.Range("E15").Value = "  +8.70%  "
$ws.Range("D16").Value = "3.046.91"
$ws.Range("E16").Value = "  +6.02%  "
$ws.Range("D17").Value = "26.13"
$ws.Range("E17").Value = "  +13.18%  "
$ws.Range("D18").Value = "2.592.49"
$ws.Range("E18").Value = "  +6.84%  "
$ws.Range("D19").Value = "8.98"
$ws.Range("E19").Value = "  +30.68%  "
$ws.Range("D20").Value = "11.75"
$ws.Range("E20").Value = "  +12.16%  "
$ws.Range("D21").Value = "375.69"
$ws.Range("E21").Value = "  +12.03%  "
$ws.Range("E22").Value = "  +19.75%  "
$ws.Range("E23").Value = "  +6.06%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "69.70"
$ws.Range("E25").Value = "  +4.38%  "
$ws.Range("E26").Value = "  +12.75%  "
$ws.Range("D27").Value = "9.21"
$ws.Range("E27").Value = "  +12.23%  "
$ws.Range("E28").Value = "  +6.40%  "
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  +14.77%  "
$ws.Range("E31").Value = "  +11.52%  "
$ws.Range("D32").Value = "506.94"
$ws.Range("E32").Value = "  +18.60%  "
$ws.Range("E33").Value = "  +18.38%  "
$ws.Range("D34").Value = "1.72"
$ws.Range("E34").Value = "  +6.94%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +12.19%  "
$ws.Range("D37").Value = "159.87"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "19.16"
$ws.Range("E38").Value = "  +7.11%  "
$ws.Range("E39").Value = "  +1.77%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  +13.19%  "
$ws.Range("D42").Value = "1.67"
$ws.Range("E42").Value = "  +12.36%  "
$ws.Range("E43").Value = "  +7.64%  "
$ws.Range("E44").Value = "  +19.39%  "
$ws.Range("D45").Value = "38.93"
$ws.Range("E45").Value = "  +4.38%  "
$ws.Range("D46").Value = "1.16"
$ws.Range("E46").Value = "  +8.17%  "
$ws.Range("D47").Value = "150.75"
$ws.Range("E47").Value = "  +14.77%  "
$ws.Range("E48").Value = "  +14.38%  "
$ws.Range("E49").Value = "  +8.17%  "
$ws.Range("E51").Value = "  +5.97%  "
